$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.034777666666667
$ws.Cells.Item(2, 8).Value = 9.104333
$ws.Cells.Item(2, 9).Value = 0.2502264227183869
$ws.Cells.Item(2, 10).Value = 0.2502264227183869
$ws.Cells.Item(2, 13).Value = 5.740110333333334
$ws.Cells.Item(2, 14).Value = 17.220331
$ws.Cells.Item(2, 15).Value = 0.2861925343043439
$ws.Cells.Item(2, 16).Value = 0.2861925343043439
$ws.Cells.Item(2, 17).Value = 17.41995864380256
$ws.Cells.Item(2, 18).Value = 156.779627794223
$ws.Cells.Item(2, 19).Value = 0.07161293406768521
$ws.Cells.Item(2, 20).Value = 0.0716129340676852

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.034777666666667
$ws.Cells.Item(3, 8).Value = 9.104333
$ws.Cells.Item(3, 9).Value = 0.2502264227183869
$ws.Cells.Item(3, 10).Value = 0.2502264227183869
$ws.Cells.Item(3, 15).Value = 0.2917347240316885
$ws.Cells.Item(3, 16).Value = 0.2917347240316885
$ws.Cells.Item(3, 17).Value = 17.757300482859
$ws.Cells.Item(3, 18).Value = 159.815704345731
$ws.Cells.Item(3, 19).Value = 0.07299973637718525
$ws.Cells.Item(3, 20).Value = 0.07299973637718525

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.034777666666667
$ws.Cells.Item(4, 8).Value = 9.104333
$ws.Cells.Item(4, 9).Value = 0.2502264227183869
$ws.Cells.Item(4, 10).Value = 0.2502264227183869
$ws.Cells.Item(4, 13).Value = 6.759986
$ws.Cells.Item(4, 14).Value = 20.279958
$ws.Cells.Item(4, 15).Value = 0.3370418707750538
$ws.Cells.Item(4, 16).Value = 0.3370418707750538
$ws.Cells.Item(4, 17).Value = 20.51505453977934
$ws.Cells.Item(4, 18).Value = 184.635490858014
$ws.Cells.Item(4, 19).Value = 0.08433678163035455
$ws.Cells.Item(4, 20).Value = 0.08433678163035455

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.034777666666667
$ws.Cells.Item(5, 8).Value = 9.104333
$ws.Cells.Item(5, 9).Value = 0.2502264227183869
$ws.Cells.Item(5, 10).Value = 0.2502264227183869
$ws.Cells.Item(5, 13).Value = 1.705448333333333
$ws.Cells.Item(5, 14).Value = 5.116345
$ws.Cells.Item(5, 15).Value = 0.0850308708889137
$ws.Cells.Item(5, 16).Value = 0.0850308708889137
$ws.Cells.Item(5, 17).Value = 5.175656513653889
$ws.Cells.Item(5, 18).Value = 46.580908622885
$ws.Cells.Item(5, 19).Value = 0.0212769706431619
$ws.Cells.Item(5, 20).Value = 0.0212769706431619

# Row 6
$ws.Cells.Item(6, 9).Value = 0.4835045831069426
$ws.Cells.Item(6, 10).Value = 0.4835045831069426
$ws.Cells.Item(6, 13).Value = 5.740110333333334
$ws.Cells.Item(6, 14).Value = 17.220331
$ws.Cells.Item(6, 15).Value = 0.2861925343043439
$ws.Cells.Item(6, 16).Value = 0.2861925343043439
$ws.Cells.Item(6, 17).Value = 33.66003378184822
$ws.Cells.Item(6, 18).Value = 302.940304036634
$ws.Cells.Item(6, 19).Value = 0.1383754019871412
$ws.Cells.Item(6, 20).Value = 0.1383754019871412

# Row 7
$ws.Cells.Item(7, 9).Value = 0.4835045831069426
$ws.Cells.Item(7, 10).Value = 0.4835045831069426
$ws.Cells.Item(7, 15).Value = 0.2917347240316885
$ws.Cells.Item(7, 16).Value = 0.2917347240316885
$ws.Cells.Item(7, 19).Value = 0.1410550761207605
$ws.Cells.Item(7, 20).Value = 0.1410550761207605

# Row 8
$ws.Cells.Item(8, 9).Value = 0.4835045831069426
$ws.Cells.Item(8, 10).Value = 0.4835045831069426
$ws.Cells.Item(8, 13).Value = 6.759986
$ws.Cells.Item(8, 14).Value = 20.279958
$ws.Cells.Item(8, 15).Value = 0.3370418707750538
$ws.Cells.Item(8, 16).Value = 0.3370418707750538
$ws.Cells.Item(8, 17).Value = 39.64058945060133
$ws.Cells.Item(8, 18).Value = 356.765305055412
$ws.Cells.Item(8, 19).Value = 0.1629612892186764
$ws.Cells.Item(8, 20).Value = 0.1629612892186764

# Row 9
$ws.Cells.Item(9, 9).Value = 0.4835045831069426
$ws.Cells.Item(9, 10).Value = 0.4835045831069426
$ws.Cells.Item(9, 13).Value = 1.705448333333333
$ws.Cells.Item(9, 14).Value = 5.116345
$ws.Cells.Item(9, 15).Value = 0.0850308708889137
$ws.Cells.Item(9, 16).Value = 0.0850308708889137
$ws.Cells.Item(9, 17).Value = 10.00075698542555
$ws.Cells.Item(9, 18).Value = 90.00681286883
$ws.Cells.Item(9, 19).Value = 0.04111281578036448
$ws.Cells.Item(9, 20).Value = 0.04111281578036448

# Row 10
$ws.Cells.Item(10, 7).Value = 2.564975
$ws.Cells.Item(10, 8).Value = 7.694925
$ws.Cells.Item(10, 9).Value = 0.2114897989601526
$ws.Cells.Item(10, 10).Value = 0.2114897989601526
$ws.Cells.Item(10, 13).Value = 5.740110333333334
$ws.Cells.Item(10, 14).Value = 17.220331
$ws.Cells.Item(10, 15).Value = 0.2861925343043439
$ws.Cells.Item(10, 16).Value = 0.2861925343043439
$ws.Cells.Item(10, 17).Value = 14.72323950224167
$ws.Cells.Item(10, 18).Value = 132.509155520175
$ws.Cells.Item(10, 19).Value = 0.06052680154392229
$ws.Cells.Item(10, 20).Value = 0.06052680154392227

# Row 11
$ws.Cells.Item(11, 7).Value = 2.564975
$ws.Cells.Item(11, 8).Value = 7.694925
$ws.Cells.Item(11, 9).Value = 0.2114897989601526
$ws.Cells.Item(11, 10).Value = 0.2114897989601526
$ws.Cells.Item(11, 15).Value = 0.2917347240316885
$ws.Cells.Item(11, 16).Value = 0.2917347240316885
$ws.Cells.Item(11, 17).Value = 15.008358703275
$ws.Cells.Item(11, 18).Value = 135.075228329475
$ws.Cells.Item(11, 19).Value = 0.06169891813515742
$ws.Cells.Item(11, 20).Value = 0.06169891813515742

# Row 12
$ws.Cells.Item(12, 7).Value = 2.564975
$ws.Cells.Item(12, 8).Value = 7.694925
$ws.Cells.Item(12, 9).Value = 0.2114897989601526
$ws.Cells.Item(12, 10).Value = 0.2114897989601526
$ws.Cells.Item(12, 13).Value = 6.759986
$ws.Cells.Item(12, 14).Value = 20.279958
$ws.Cells.Item(12, 15).Value = 0.3370418707750538
$ws.Cells.Item(12, 16).Value = 0.3370418707750538
$ws.Cells.Item(12, 17).Value = 17.33919509035
$ws.Cells.Item(12, 18).Value = 156.05275581315
$ws.Cells.Item(12, 19).Value = 0.07128091749136987
$ws.Cells.Item(12, 20).Value = 0.07128091749136987

# Row 13
$ws.Cells.Item(13, 7).Value = 2.564975
$ws.Cells.Item(13, 8).Value = 7.694925
$ws.Cells.Item(13, 9).Value = 0.2114897989601526
$ws.Cells.Item(13, 10).Value = 0.2114897989601526
$ws.Cells.Item(13, 13).Value = 1.705448333333333
$ws.Cells.Item(13, 14).Value = 5.116345
$ws.Cells.Item(13, 15).Value = 0.0850308708889137
$ws.Cells.Item(13, 16).Value = 0.0850308708889137
$ws.Cells.Item(13, 17).Value = 4.374432338791666
$ws.Cells.Item(13, 18).Value = 39.36989104912499
$ws.Cells.Item(13, 19).Value = 0.01798316178970305
$ws.Cells.Item(13, 20).Value = 0.01798316178970305

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.6643690000000001
$ws.Cells.Item(14, 8).Value = 1.993107
$ws.Cells.Item(14, 9).Value = 0.05477919521451775
$ws.Cells.Item(14, 10).Value = 0.05477919521451775
$ws.Cells.Item(14, 13).Value = 5.740110333333334
$ws.Cells.Item(14, 14).Value = 17.220331
$ws.Cells.Item(14, 15).Value = 0.2861925343043439
$ws.Cells.Item(14, 16).Value = 0.2861925343043439
$ws.Cells.Item(14, 17).Value = 3.813551362046334
$ws.Cells.Item(14, 18).Value = 34.32196225841701
$ws.Cells.Item(14, 19).Value = 0.01567739670559522
$ws.Cells.Item(14, 20).Value = 0.01567739670559522

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.6643690000000001
$ws.Cells.Item(15, 8).Value = 1.993107
$ws.Cells.Item(15, 9).Value = 0.05477919521451775
$ws.Cells.Item(15, 10).Value = 0.05477919521451775
$ws.Cells.Item(15, 15).Value = 0.2917347240316885
$ws.Cells.Item(15, 16).Value = 0.2917347240316885
$ws.Cells.Item(15, 17).Value = 3.887401734261
$ws.Cells.Item(15, 18).Value = 34.986615608349
$ws.Cells.Item(15, 19).Value = 0.01598099339858533
$ws.Cells.Item(15, 20).Value = 0.01598099339858533

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.6643690000000001
$ws.Cells.Item(16, 8).Value = 1.993107
$ws.Cells.Item(16, 9).Value = 0.05477919521451775
$ws.Cells.Item(16, 10).Value = 0.05477919521451775
$ws.Cells.Item(16, 13).Value = 6.759986
$ws.Cells.Item(16, 14).Value = 20.279958
$ws.Cells.Item(16, 15).Value = 0.3370418707750538
$ws.Cells.Item(16, 16).Value = 0.3370418707750538
$ws.Cells.Item(16, 17).Value = 4.491125138834001
$ws.Cells.Item(16, 18).Value = 40.420126249506
$ws.Cells.Item(16, 19).Value = 0.01846288243465294
$ws.Cells.Item(16, 20).Value = 0.01846288243465294

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.6643690000000001
$ws.Cells.Item(17, 8).Value = 1.993107
$ws.Cells.Item(17, 9).Value = 0.05477919521451775
$ws.Cells.Item(17, 10).Value = 0.05477919521451775
$ws.Cells.Item(17, 13).Value = 1.705448333333333
$ws.Cells.Item(17, 14).Value = 5.116345
$ws.Cells.Item(17, 15).Value = 0.0850308708889137
$ws.Cells.Item(17, 16).Value = 0.0850308708889137
$ws.Cells.Item(17, 17).Value = 1.133047003768333
$ws.Cells.Item(17, 18).Value = 10.197423033915
$ws.Cells.Item(17, 19).Value = 0.004657922675684259
$ws.Cells.Item(17, 20).Value = 0.004657922675684258
